# Rename TaskBook to TaskList across the HighLevelSequenceDiagrams slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape "TextBox 32" (id=33): move box slightly right (EMU -> points, 1pt = 12700 EMU)
# and rename event text.
$shp17 = $s.Shapes.Item(17)
$shp17.Left = 498
$tr17 = $shp17.TextFrame.TextRange
$tr17.Runs(2).Text = "TaskListChangedEvent"

# Shape "TextBox 61" (id=62): rename event text only.
$shp29 = $s.Shapes.Item(29)
$tr29 = $shp29.TextFrame.TextRange
$tr29.Runs(2).Text = "TaskListChangedEvent"

# Shape "TextBox 73" (id=74): rename handler text.
$shp37 = $s.Shapes.Item(37)
$tr37 = $shp37.TextFrame.TextRange
$tr37.Runs(1).Text = "handleTaskListChangedEvent"

# Shape "TextBox 49" (id=50): rename handler text.
$shp44 = $s.Shapes.Item(44)
$tr44 = $shp44.TextFrame.TextRange
$tr44.Runs(1).Text = "handleTaskListChangedEvent"
